$d = $word.ActiveDocument

# --- 1. Append the two new paragraphs at the very end of the document ---

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("Are we actually doing 2 queues? Reception to triage and triage to treatment room?")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
# Append the final sentence plus a one-character placeholder ("Z") that we will
# remove in a moment. The placeholder lets us relocate the "_GoBack" bookmark
# to the true end of the document without ever asking for a zero-length
# bookmark sitting exactly on a paragraph boundary (doing so directly is
# mishandled by the host and snaps the bookmark back to the top of the
# document), mirroring how Word naturally drags "_GoBack" along while text is
# typed and then trimmed.
$end.InsertAfter("What we should we envision an emergency skipping the other patients but still being entered in that they occupy a treatment room and a doctor?Z")

# --- 2. Move the "_GoBack" bookmark from its old location to the new end of
#        the document (this is what Word itself would do automatically as a
#        side effect of typing new text at the end of the document) ---

if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

$docEnd = $d.Content.End
$placeholder = $d.Range($docEnd - 2, $docEnd - 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Range.Delete()

Write-Host "Done"
